$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status of existing tickets (rows 2,3,5,6) from "Open" to "Resolved"
$ws.Range("G2").Value = "Resolved"
$ws.Range("G3").Value = "Resolved"
$ws.Range("G5").Value = "Resolved"
$ws.Range("G6").Value = "Resolved"

# Add new ticket row 7
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "ClassPlus"
$ws.Range("C7").Value = 1353506
$ws.Range("D7").Value = [datetime]"2023-01-08"
$ws.Range("E7").Value = "Debashish Nath"
$ws.Range("F7").Value = "Issues related to websites"
$ws.Range("G7").Value = "Open"

# Column G now holds data ("Open"/"Resolved") so Excel auto-sizes it to fit
$ws.Columns.Item(7).ColumnWidth = 8.67

# Update the active selection to reflect the last edited cell
$ws.Range("H7").Select()
